$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 3.85
$ws.Range("H7").Value = 3.25
$ws.Range("I7").Value = 1.88
$ws.Range("J7").Value = 4.3
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 2.52
$ws.Range("N7").Value = 6.7
$ws.Range("O7").Value = 1.35
$ws.Range("P7").Value = 2.95
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.67
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.85
$ws.Range("W7").Value = 10.5
$ws.Range("X7").Value = 21
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 60
$ws.Range("AA7").Value = 37
$ws.Range("AC7").Value = 6.7
$ws.Range("AD7").Value = 6.4
$ws.Range("AF7").Value = 75
$ws.Range("AH7").Value = 8.5
$ws.Range("AJ7").Value = 16
$ws.Range("AK7").Value = 16
$ws.Range("AM7").Value = 600
$ws.Range("AN7").Value = 5.7
$ws.Range("AO7").Value = 22
$ws.Range("AT7").Value = 2.67
$ws.Range("AU7").Value = 7.2
$ws.Range("AW7").Value = 3.75
$ws.Range("AX7").Value = 9.75
$ws.Range("AZ7").Value = 37
$ws.Range("BA7").Value = 75
